$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (prices with "." as thousands/decimal
# separators). Force Text number format before assignment so Excel does not
# auto-convert the string to a number (which would drop trailing zeros /
# reformat the separators), then restore the cell to the default "Normal"
# style so no stray style index lingers on the cell.
function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "60.335.65"
$ws.Range("E2").Value = "  +2.27%  "
Set-TextValue $ws "D3" "2.595.43"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("E4").Value = "  +0.07%  "
Set-TextValue $ws "D5" "583.60"
$ws.Range("E5").Value = "  +6.57%  "
Set-TextValue $ws "D6" "142.83"
$ws.Range("E6").Value = "  +2.41%  "
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("E8").Value = "  +1.74%  "
Set-TextValue $ws "D9" "2.615.72"
$ws.Range("E9").Value = "  +2.73%  "
Set-TextValue $ws "D10" "6.53"
$ws.Range("E10").Value = "  -1.88%  "
Set-TextValue $ws "D11" "0.105"
$ws.Range("E11").Value = "  +2.89%  "
$ws.Range("E12").Value = "  -3.11%  "
Set-TextValue $ws "D13" "0.371"
$ws.Range("E13").Value = "  +5.89%  "
Set-TextValue $ws "D14" "3.071.61"
$ws.Range("E14").Value = "  +2.53%  "
Set-TextValue $ws "D15" "24.53"
$ws.Range("E15").Value = "  +6.82%  "
Set-TextValue $ws "D16" "60.382.79"
$ws.Range("E16").Value = "  +2.30%  "
$ws.Range("E17").Value = "  +4.21%  "
Set-TextValue $ws "D18" "2.619.01"
$ws.Range("E18").Value = "  +2.23%  "
Set-TextValue $ws "D19" "11.34"
$ws.Range("E19").Value = "  +11.31%  "
Set-TextValue $ws "D20" "4.67"
$ws.Range("E20").Value = "  +3.62%  "
Set-TextValue $ws "D21" "347.74"
$ws.Range("E21").Value = "  +3.81%  "
$ws.Range("E22").Value = "  +8.59%  "
$ws.Range("E23").Value = "  +0.83%  "
Set-TextValue $ws "D24" "0.521"
$ws.Range("E24").Value = "  +10.06%  "
Set-TextValue $ws "D25" "63.05"
$ws.Range("E25").Value = "  +1.08%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("E27").Value = "  +1.43%  "
Set-TextValue $ws "D28" "7.93"
$ws.Range("E28").Value = "  +8.12%  "
Set-TextValue $ws "D29" "0.0₃0796"
$ws.Range("E29").Value = "  +5.48%  "
Set-TextValue $ws "D30" "1.86"
$ws.Range("E30").Value = "  +12.17%  "
$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws "D31" "0.998"
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws "D32" "6.37"
$ws.Range("E32").Value = "  +4.36%  "
Set-TextValue $ws "D33" "164.29"
$ws.Range("E33").Value = "  +3.92%  "
Set-TextValue $ws "D34" "19.50"
$ws.Range("E34").Value = "  +3.26%  "
Set-TextValue $ws "D35" "4.24"
$ws.Range("E35").Value = "  +5.19%  "
Set-TextValue $ws "D36" "0.990"
$ws.Range("E36").Value = "  +12.26%  "
$ws.Range("E37").Value = "  +7.25%  "
$ws.Range("E38").Value = "  +11.26%  "
Set-TextValue $ws "D39" "37.94"
$ws.Range("E39").Value = "  +1.80%  "
Set-TextValue $ws "D40" "311.19"
$ws.Range("E40").Value = "  +10.41%  "
$ws.Range("E41").Value = "  +7.17%  "
$ws.Range("E42").Value = "  +0.52%  "
Set-TextValue $ws "D43" "135.66"
$ws.Range("E43").Value = "  +0.27%  "
Set-TextValue $ws "D44" "5.08"
$ws.Range("E44").Value = "  +14.82%  "
Set-TextValue $ws "D45" "0.0990"
$ws.Range("E45").Value = "  +2.47%  "
Set-TextValue $ws "D46" "0.997"
$ws.Range("E46").Value = "  -0.22%  "
Set-TextValue $ws "D47" "19.82"
$ws.Range("E47").Value = "  +6.93%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws "D48" "0.0551"
$ws.Range("E48").Value = "  +5.05%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws "D49" "0.604"
$ws.Range("E49").Value = "  +3.44%  "
Set-TextValue $ws "D50" "20.17"
$ws.Range("E50").Value = "  +9.86%  "
Set-TextValue $ws "D51" "0.0241"
$ws.Range("E51").Value = "  +4.25%  "
